$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Init")

# The "Lower Right Cell" column (D) references shift down one row,
# from row 25 to row 26, for the block of indices used by the
# delivery-data extraction (scenario 46 init update).
$ws.Range("D5").Value  = "A26"
$ws.Range("D6").Value  = "B26"
$ws.Range("D7").Value  = "C26"
$ws.Range("D8").Value  = "G26"
$ws.Range("D9").Value  = "H26"
$ws.Range("D10").Value = "I26"
$ws.Range("D11").Value = "J26"
